$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# Remove the "Employee ID- " paragraph entirely. Deleting the paragraph's
# whole Range (text + trailing paragraph mark) merges it with the next
# paragraph ("D/o Sh. O. P. Shivhare"), which is exactly what the diff
# shows: the "Employee ID- " run and its paragraph break disappear.
$empIdPara = $d.Paragraphs.Item(7)
$empIdPara.Range.Delete()

# --- Change 2 -----------------------------------------------------------
# After the "Your contributions ..." paragraph there were two blank
# paragraphs before "Thanking you,", and the _GoBack bookmark sat at the
# end of the "Your contributions ..." paragraph. The edit inserts a new
# blank paragraph right after "Your contributions ..." and moves the
# bookmark to what is now the following (first pre-existing) blank
# paragraph.
$contribPara = $d.Paragraphs.Item(12)

# Insert a clean blank paragraph (matching the sz=24 paragraph mark
# formatting used elsewhere) right after the "Your contributions ..."
# paragraph, without leaving a stray empty run behind.
$insertionPoint = $d.Range($contribPara.Range.End, $contribPara.Range.End)
$newParaXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newParaXml)

# Move the _GoBack bookmark off the "Your contributions ..." paragraph
# and onto the next (now second) blank paragraph.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$targetPara = $d.Paragraphs.Item(14)
$targetRange = $targetPara.Range.Duplicate()
$targetRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $targetRange)
